# Gestion de Cambios - add new change-request entry (row 10) and
# backfill the "Firma" (M) column with "PMO" for the existing entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108
$xlLeft   = -4131

# --- 1. Backfill column M ("Firma") for rows 3-9 with "PMO" --------------
# These rows already contain data; only the Firma cell was empty. Excel
# centers that column (horizontal + vertical), matching the other filled
# cells in the same row.
foreach ($r in 3..9) {
    $cell = $ws.Cells.Item($r, 13)   # column M
    $cell.Value = "PMO"
    $cell.HorizontalAlignment = $xlCenter
    $cell.VerticalAlignment = $xlCenter
}

# Row 9's "Coste" (J) and "Aprobado por" (L) cells also become centered
# (horizontal+vertical) along with the rest of the now-complete row.
$ws.Cells.Item(9, 10).HorizontalAlignment = $xlCenter   # J9
$ws.Cells.Item(9, 10).VerticalAlignment = $xlCenter
$ws.Cells.Item(9, 12).HorizontalAlignment = $xlCenter   # L9
$ws.Cells.Item(9, 12).VerticalAlignment = $xlCenter

# --- 2. Fill in the new entry on row 10 -----------------------------------
# C10 needs the same date format as the other "Fecha actual" cells above it
# (numFmtId 14, centered). Copy the format from C9 rather than re-stating a
# NumberFormat string, so the existing style is reused instead of a new one
# being minted.
$ws.Cells.Item(9, 3).Copy()
$ws.Cells.Item(10, 3).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(10, 3).Value = 43428                       # C10 Fecha actual

$ws.Cells.Item(10, 4).Value = 1                            # D10 Version ant
$ws.Cells.Item(10, 4).HorizontalAlignment = $xlCenter
$ws.Cells.Item(10, 4).VerticalAlignment = $xlCenter

$ws.Cells.Item(10, 5).Value = 1.1000000000000001           # E10 Version nueva
$ws.Cells.Item(10, 5).HorizontalAlignment = $xlCenter
$ws.Cells.Item(10, 5).VerticalAlignment = $xlCenter

$ws.Cells.Item(10, 6).Value = "PMO"                         # F10 Solicitador
$ws.Cells.Item(10, 6).HorizontalAlignment = $xlCenter
$ws.Cells.Item(10, 6).VerticalAlignment = $xlCenter

$ws.Cells.Item(10, 7).Value = "Correción y adición de mapas de procesos a E_3"   # G10 Descripcion
$ws.Cells.Item(10, 8).Value = "Anexo RRHH1 y 2"                                  # H10 Paquete WBS afectado
$ws.Cells.Item(10, 9).Value = "Nuevos mapas de procesos, corrección en el uso de referencias"  # I10 Justificacion

$ws.Cells.Item(10, 10).Value = "-"                          # J10 Coste
$ws.Cells.Item(10, 10).HorizontalAlignment = $xlCenter
$ws.Cells.Item(10, 10).VerticalAlignment = $xlCenter

$ws.Cells.Item(10, 12).Value = "PMO"                        # L10 Aprobado por
$ws.Cells.Item(10, 12).HorizontalAlignment = $xlCenter
$ws.Cells.Item(10, 12).VerticalAlignment = $xlCenter

$ws.Cells.Item(10, 13).Value = "PMO"                        # M10 Firma
$ws.Cells.Item(10, 13).HorizontalAlignment = $xlCenter
$ws.Cells.Item(10, 13).VerticalAlignment = $xlCenter

# --- 3. Move the active selection to D10:E10 ------------------------------
$ws.Range("D10:E10").Select()
